$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7317783236503601
$ws.Range("B1").Value = 1.142157316207886
$ws.Range("C1").Value = 2.477017641067505
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.815943479537964
